# Natmi following Dr Hou advice
# Rewrite the LR-pair data table: a new "ECs" sending/target cluster group is
# introduced (string table gains "ECs" right after the header strings, ahead of
# "FAPs"), and the Sending/Target-cluster x Ligand/Receptor combinations are
# expanded from 2x3 to 3x3 rows (rows 2-10) with refreshed metric values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Ntrk3"
$ws.Cells.Item(2, 3).Value = "Ptprf"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.138349
$ws.Cells.Item(2, 8).Value = 0.415047
$ws.Cells.Item(2, 9).Value = 0.04449686021752534
$ws.Cells.Item(2, 10).Value = 0.04449686021752534
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.21127
$ws.Cells.Item(2, 14).Value = 0.63381
$ws.Cells.Item(2, 15).Value = 0.02468048274231428
$ws.Cells.Item(2, 16).Value = 0.02468048274231428
$ws.Cells.Item(2, 17).Value = 0.02922899323
$ws.Cells.Item(2, 18).Value = 0.26306093907
$ws.Cells.Item(2, 19).Value = 0.001098203990685805
$ws.Cells.Item(2, 20).Value = 0.001098203990685805

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Ntrk3"
$ws.Cells.Item(3, 3).Value = "Ptprf"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.138349
$ws.Cells.Item(3, 8).Value = 0.415047
$ws.Cells.Item(3, 9).Value = 0.04449686021752534
$ws.Cells.Item(3, 10).Value = 0.04449686021752534
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 3.467027333333334
$ws.Cells.Item(3, 14).Value = 10.401082
$ws.Cells.Item(3, 15).Value = 0.405016842275123
$ws.Cells.Item(3, 16).Value = 0.405016842275123
$ws.Cells.Item(3, 17).Value = 0.4796597645393334
$ws.Cells.Item(3, 18).Value = 4.316937880854001
$ws.Cells.Item(3, 19).Value = 0.01802197781645966
$ws.Cells.Item(3, 20).Value = 0.01802197781645965

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Ntrk3"
$ws.Cells.Item(4, 3).Value = "Ptprf"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.138349
$ws.Cells.Item(4, 8).Value = 0.415047
$ws.Cells.Item(4, 9).Value = 0.04449686021752534
$ws.Cells.Item(4, 10).Value = 0.04449686021752534
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 4.881908
$ws.Cells.Item(4, 14).Value = 14.645724
$ws.Cells.Item(4, 15).Value = 0.5703026749825627
$ws.Cells.Item(4, 16).Value = 0.5703026749825627
$ws.Cells.Item(4, 17).Value = 0.675407089892
$ws.Cells.Item(4, 18).Value = 6.078663809028001
$ws.Cells.Item(4, 19).Value = 0.02537667841037988
$ws.Cells.Item(4, 20).Value = 0.02537667841037988

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Ntrk3"
$ws.Cells.Item(5, 3).Value = "Ptprf"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 2.468673666666667
$ws.Cells.Item(5, 8).Value = 7.406021
$ws.Cells.Item(5, 9).Value = 0.7939936469967431
$ws.Cells.Item(5, 10).Value = 0.793993646996743
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.21127
$ws.Cells.Item(5, 14).Value = 0.63381
$ws.Cells.Item(5, 15).Value = 0.02468048274231428
$ws.Cells.Item(5, 16).Value = 0.02468048274231428
$ws.Cells.Item(5, 17).Value = 0.5215566855566667
$ws.Cells.Item(5, 18).Value = 4.694010170009999
$ws.Cells.Item(5, 19).Value = 0.01959614650221029
$ws.Cells.Item(5, 20).Value = 0.01959614650221029

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Ntrk3"
$ws.Cells.Item(6, 3).Value = "Ptprf"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 2.468673666666667
$ws.Cells.Item(6, 8).Value = 7.406021
$ws.Cells.Item(6, 9).Value = 0.7939936469967431
$ws.Cells.Item(6, 10).Value = 0.793993646996743
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 3.467027333333334
$ws.Cells.Item(6, 14).Value = 10.401082
$ws.Cells.Item(6, 15).Value = 0.405016842275123
$ws.Cells.Item(6, 16).Value = 0.405016842275123
$ws.Cells.Item(6, 17).Value = 8.558959079413556
$ws.Cells.Item(6, 18).Value = 77.030631714722
$ws.Cells.Item(6, 19).Value = 0.3215807996931296
$ws.Cells.Item(6, 20).Value = 0.3215807996931296

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Ntrk3"
$ws.Cells.Item(7, 3).Value = "Ptprf"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 2.468673666666667
$ws.Cells.Item(7, 8).Value = 7.406021
$ws.Cells.Item(7, 9).Value = 0.7939936469967431
$ws.Cells.Item(7, 10).Value = 0.793993646996743
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 4.881908
$ws.Cells.Item(7, 14).Value = 14.645724
$ws.Cells.Item(7, 15).Value = 0.5703026749825627
$ws.Cells.Item(7, 16).Value = 0.5703026749825627
$ws.Cells.Item(7, 17).Value = 12.05183772268933
$ws.Cells.Item(7, 18).Value = 108.466539504204
$ws.Cells.Item(7, 19).Value = 0.4528167008014032
$ws.Cells.Item(7, 20).Value = 0.4528167008014032

# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Ntrk3"
$ws.Cells.Item(8, 3).Value = "Ptprf"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.5021629999999999
$ws.Cells.Item(8, 8).Value = 1.506489
$ws.Cells.Item(8, 9).Value = 0.1615094927857316
$ws.Cells.Item(8, 10).Value = 0.1615094927857315
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 0.21127
$ws.Cells.Item(8, 14).Value = 0.63381
$ws.Cells.Item(8, 15).Value = 0.02468048274231428
$ws.Cells.Item(8, 16).Value = 0.02468048274231428
$ws.Cells.Item(8, 17).Value = 0.10609197701
$ws.Cells.Item(8, 18).Value = 0.9548277930899998
$ws.Cells.Item(8, 19).Value = 0.00398613224941818
$ws.Cells.Item(8, 20).Value = 0.00398613224941818

# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Ntrk3"
$ws.Cells.Item(9, 3).Value = "Ptprf"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.5021629999999999
$ws.Cells.Item(9, 8).Value = 1.506489
$ws.Cells.Item(9, 9).Value = 0.1615094927857316
$ws.Cells.Item(9, 10).Value = 0.1615094927857315
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 3.467027333333334
$ws.Cells.Item(9, 14).Value = 10.401082
$ws.Cells.Item(9, 15).Value = 0.405016842275123
$ws.Cells.Item(9, 16).Value = 0.405016842275123
$ws.Cells.Item(9, 17).Value = 1.741012846788667
$ws.Cells.Item(9, 18).Value = 15.669115621098
$ws.Cells.Item(9, 19).Value = 0.06541406476553376
$ws.Cells.Item(9, 20).Value = 0.06541406476553376

# Row 10
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Ntrk3"
$ws.Cells.Item(10, 3).Value = "Ptprf"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.5021629999999999
$ws.Cells.Item(10, 8).Value = 1.506489
$ws.Cells.Item(10, 9).Value = 0.1615094927857316
$ws.Cells.Item(10, 10).Value = 0.1615094927857315
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 4.881908
$ws.Cells.Item(10, 14).Value = 14.645724
$ws.Cells.Item(10, 15).Value = 0.5703026749825627
$ws.Cells.Item(10, 16).Value = 0.5703026749825627
$ws.Cells.Item(10, 17).Value = 2.451513567004
$ws.Cells.Item(10, 18).Value = 22.063622103036
$ws.Cells.Item(10, 19).Value = 0.09210929577077963
$ws.Cells.Item(10, 20).Value = 0.09210929577077961
